# Recalculate "Median Value" (score relative to median) and "Tier" columns
# after merging with zip/census tract data. The underlying median used to
# normalize each school's raw value changed, so every ratio in column C is
# rescaled, and the quartile-based Tier labels in column D are reassigned
# based on the new ranking of values that remain at/above the median.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; MedianValue = 1.135083983912941; Tier = "3rd Tier" },
    @{ Row = 3; MedianValue = 0.6841731724627396; Tier = "Below Median" },
    @{ Row = 4; MedianValue = 1.685592618878637; Tier = "1st Tier" },
    @{ Row = 5; MedianValue = 1.376389874615566; Tier = "2nd Tier" },
    @{ Row = 6; MedianValue = 1.685592618878637; Tier = "1st Tier" },
    @{ Row = 7; MedianValue = 2.045422285308729; Tier = "1st Tier" },
    @{ Row = 8; MedianValue = 1.325762952448545; Tier = "2nd Tier" },
    @{ Row = 9; MedianValue = 2.045422285308729; Tier = "1st Tier" },
    @{ Row = 10; MedianValue = 0.3553347527797492; Tier = "Below Median" },
    @{ Row = 11; MedianValue = 0.4178850248403123; Tier = "Below Median" },
    @{ Row = 12; MedianValue = 0.6728176011355571; Tier = "Below Median" },
    @{ Row = 13; MedianValue = 0.9609652235628106; Tier = "Below Median" },
    @{ Row = 14; MedianValue = 1.163946061036196; Tier = "3rd Tier" },
    @{ Row = 15; MedianValue = 1.163946061036196; Tier = "3rd Tier" },
    @{ Row = 16; MedianValue = 0.7324343506032647; Tier = "Below Median" },
    @{ Row = 17; MedianValue = 0.9247693399574166; Tier = "Below Median" },
    @{ Row = 18; MedianValue = 0.7686302342086586; Tier = "Below Median" },
    @{ Row = 19; MedianValue = 0.8833688194937308; Tier = "Below Median" },
    @{ Row = 20; MedianValue = 0.8833688194937308; Tier = "Below Median" },
    @{ Row = 21; MedianValue = 1.185237757274663; Tier = "3rd Tier" },
    @{ Row = 22; MedianValue = 0.5613910574875799; Tier = "Below Median" },
    @{ Row = 23; MedianValue = 0.5613910574875799; Tier = "Below Median" },
    @{ Row = 24; MedianValue = 1.182398864442867; Tier = "3rd Tier" },
    @{ Row = 25; MedianValue = 1.611071682044003; Tier = "2nd Tier" },
    @{ Row = 26; MedianValue = 1.753726046841732; Tier = "1st Tier" },
    @{ Row = 27; MedianValue = 1.753726046841732; Tier = "1st Tier" },
    @{ Row = 28; MedianValue = 1.474095102909865; Tier = "2nd Tier" },
    @{ Row = 29; MedianValue = 0.8211497515968772; Tier = "Below Median" },
    @{ Row = 30; MedianValue = 1.249112845990064; Tier = "2nd Tier" },
    @{ Row = 31; MedianValue = 0.9141234918381831; Tier = "Below Median" },
    @{ Row = 32; MedianValue = 1.182398864442867; Tier = "3rd Tier" },
    @{ Row = 33; MedianValue = 1.039744499645138; Tier = "4th Tier" },
    @{ Row = 34; MedianValue = 0.9921930447125621; Tier = "Below Median" },
    @{ Row = 35; MedianValue = 0.7314880529926663; Tier = "Below Median" },
    @{ Row = 36; MedianValue = 1.028152353915306; Tier = "4th Tier" },
    @{ Row = 37; MedianValue = 0.9027679205110007; Tier = "Below Median" },
    @{ Row = 38; MedianValue = 0.9027679205110007; Tier = "Below Median" },
    @{ Row = 39; MedianValue = 1.462029808374734; Tier = "2nd Tier" },
    @{ Row = 40; MedianValue = 1.028152353915306; Tier = "4th Tier" },
    @{ Row = 41; MedianValue = 0.9027679205110007; Tier = "Below Median" },
    @{ Row = 42; MedianValue = 0.9623846699787083; Tier = "Below Median" },
    @{ Row = 43; MedianValue = 1.182398864442867; Tier = "3rd Tier" },
    @{ Row = 44; MedianValue = 1.12230896616986; Tier = "4th Tier" },
    @{ Row = 45; MedianValue = 1.135083983912941; Tier = "3rd Tier" },
    @{ Row = 46; MedianValue = 0.6898509581263307; Tier = "Below Median" },
    @{ Row = 47; MedianValue = 1.039744499645138; Tier = "4th Tier" },
    @{ Row = 48; MedianValue = 0.890347764371895; Tier = "Below Median" },
    @{ Row = 49; MedianValue = 0.7885024840312278; Tier = "Below Median" },
    @{ Row = 50; MedianValue = 0.7314880529926663; Tier = "Below Median" },
    @{ Row = 51; MedianValue = 1.135083983912941; Tier = "3rd Tier" },
    @{ Row = 52; MedianValue = 1.039744499645138; Tier = "4th Tier" },
    @{ Row = 53; MedianValue = 1.090844570617459; Tier = "4th Tier" },
    @{ Row = 54; MedianValue = 0.4059616749467708; Tier = "Below Median" },
    @{ Row = 55; MedianValue = 0.9854506742370476; Tier = "Below Median" },
    @{ Row = 56; MedianValue = 1.235391530636385; Tier = "3rd Tier" },
    @{ Row = 57; MedianValue = 0.5734563520227112; Tier = "Below Median" },
    @{ Row = 58; MedianValue = 0.5956943458717767; Tier = "Below Median" },
    @{ Row = 59; MedianValue = 0.9666430092264017; Tier = "Below Median" },
    @{ Row = 60; MedianValue = 0.850958126330731; Tier = "Below Median" },
    @{ Row = 61; MedianValue = 1.325762952448545; Tier = "2nd Tier" },
    @{ Row = 62; MedianValue = 1.462029808374734; Tier = "2nd Tier" },
    @{ Row = 63; MedianValue = 0.9623846699787083; Tier = "Below Median" },
    @{ Row = 64; MedianValue = 0.5956943458717767; Tier = "Below Median" },
    @{ Row = 65; MedianValue = 1.039744499645138; Tier = "4th Tier" },
    @{ Row = 66; MedianValue = 0.8034066713981547; Tier = "Below Median" },
    @{ Row = 67; MedianValue = 1.082564466524722; Tier = "4th Tier" },
    @{ Row = 68; MedianValue = 1.611071682044003; Tier = "1st Tier" },
    @{ Row = 69; MedianValue = 0.3222143364088006; Tier = "Below Median" },
    @{ Row = 70; MedianValue = 0.6174591909155429; Tier = "Below Median" },
    @{ Row = 71; MedianValue = 1.611071682044003; Tier = "1st Tier" },
    @{ Row = 72; MedianValue = 0.9609652235628106; Tier = "Below Median" },
    @{ Row = 73; MedianValue = 0.4932576295244854; Tier = "Below Median" },
    @{ Row = 74; MedianValue = 1.323278921220724; Tier = "2nd Tier" },
    @{ Row = 75; MedianValue = 1.685592618878637; Tier = "1st Tier" },
    @{ Row = 76; MedianValue = 1.611071682044003; Tier = "1st Tier" },
    @{ Row = 77; MedianValue = 1; Tier = "4th Tier" },
    @{ Row = 78; MedianValue = 1.470073338064821; Tier = "2nd Tier" },
    @{ Row = 79; MedianValue = 0.971611071682044; Tier = "Below Median" },
    @{ Row = 80; MedianValue = 0.4758694109297374; Tier = "Below Median" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.MedianValue
    $ws.Cells.Item($u.Row, 4).Value = $u.Tier
}
